$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "O atendente clica em liberar veiculo inicia o (CdU004)."
#    veiculo -> veículo (also drops the spell-check markup, which happens
#    naturally once the misspelling is fixed)
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(46)
$p1.Range.Find.Execute("veiculo", $true, $false, $false, $false, $false, $true, 1, $false, "veículo", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "... informa ao gerente que o ano não e valido." -> "... não é valido."
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(74)
$p2.Range.Find.Execute("ano não e valido", $true, $false, $false, $false, $false, $true, 1, $false, "ano não é valido", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Permite que o gerente veículo da frota da empresa." ->
#    "Permite que o gerente exclua veículo da frota." with a _GoBack
#    bookmark left right before the final period (mirrors Word's own
#    behaviour of marking the last edit location).
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(80)
$p3.Range.Find.Execute("gerente veículo da frota da empresa.", $true, $false, $false, $false, $false, $true, 1, $false, "gerente exclua veículo da frota.", 2) | Out-Null

# Re-find the trailing period we just produced so we can drop a _GoBack
# bookmark immediately before it, matching the target markup.
$period = $p3.Range.Duplicate
$period.Find.Execute(".", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackRange = $d.Range($period.Start, $period.Start)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# ---------------------------------------------------------------------------
# 4) Drop the _GoBack bookmark that used to sit at the very end of the
#    document (after "O gerente clica em “cancelar” o Caso de uso e
#    encerrado.") -- it has effectively moved to change (3) above.
# ---------------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($lastParaIndex)
# nothing else to do here; the single _GoBack bookmark allowed in the
# document already got re-anchored to its new location in step 3.
